$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# MarketWebPortal row: add "Swagger" before "SOA"
$ws.Range("B2").Value = "C#, .NET Core, MVC, Web API, REST, ADO.NET, xUnit testing, Moq, Bootstrap, Unobtrusive JavaScript, Swagger, SOA"

# StarwarsWebPortal row: add ", Nlog" at the end
$ws.Range("B3").Value = "C#, .NET Core, MVC, Web API, REST, ADO.NET, xUnit testing, Moq, Bootstrap, Unobtrusive JavaScript, Nlog"
